$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list with latest scraped price/volume data.
# Column D (Price) values are forced to Text via a leading quote prefix,
# matching the original inlineStr text cells (e.g. "29.385.35" must stay
# text, not be parsed as a number).

$ws.Range("D2").Value = "'29.385.35"
$ws.Range("E2").Value = '  -0.31%  '

$ws.Range("D3").Value = "'1.847.03"
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("D4").Value = "'0.9988"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = "'240.47"
$ws.Range("E5").Value = '  -1.06%  '

$ws.Range("D6").Value = "'0.6331"
$ws.Range("E6").Value = '  -3.92%  '

$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = "'0.07562"
$ws.Range("E8").Value = '  +1.03%  '

$ws.Range("D9").Value = "'0.2966"
$ws.Range("E9").Value = '  -0.88%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = "'2.862.42"
$ws.Range("E10").Value = '  +54.88%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = "'24.54"
$ws.Range("E11").Value = '  +0.80%  '

$ws.Range("D12").Value = "'0.07714"
$ws.Range("E12").Value = '  +0.95%  '

$ws.Range("D13").Value = "'4.983"
$ws.Range("E13").Value = '  -0.77%  '

$ws.Range("D14").Value = "'0.6844"
$ws.Range("E14").Value = '  -0.02%  '

$ws.Range("D15").Value = "'0.000009955"
$ws.Range("E15").Value = '  +4.57%  '

$ws.Range("D16").Value = "'82.75"
$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("D17").Value = "'6.186"
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").Value = "'29.417.76"
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("D19").Value = "'231.41"
$ws.Range("E19").Value = '  -2.32%  '

$ws.Range("D20").Value = "'12.51"
$ws.Range("E20").Value = '  -0.62%  '

$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("D22").Value = "'7.585"
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Value = "'154.74"
$ws.Range("E24").Value = '  -1.32%  '

$ws.Range("D25").Value = "'0.1392"
$ws.Range("E25").Value = '  -2.26%  '

$ws.Range("E26").Value = '  -0.81%  '

$ws.Range("D27").Value = "'17.66"
$ws.Range("E27").Value = '  -0.75%  '

$ws.Range("D28").Value = "'1.470"
$ws.Range("E28").Value = '  -1.33%  '

$ws.Range("D29").Value = "'0.05809"
$ws.Range("E29").Value = '  -3.63%  '

$ws.Range("D30").Value = "'1.260"
$ws.Range("E30").Value = '  +0.85%  '

$ws.Range("D31").Value = "'4.118"
$ws.Range("E31").Value = '  -0.47%  '

$ws.Range("E32").Value = '  -1.47%  '

$ws.Range("D33").Value = "'1.867"
$ws.Range("E33").Value = '  +0.73%  '

$ws.Range("B34").Value = 'RocketPoolETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D34").Value = "'2.977.55"
$ws.Range("E34").Value = '  +48.41%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = "'1.159"
$ws.Range("E35").Value = '  -1.78%  '

$ws.Range("E36").Value = '  -0.89%  '

$ws.Range("D37").Value = "'2.593"
$ws.Range("E37").Value = '  -0.10%  '

$ws.Range("D38").Value = "'1.250.93"
$ws.Range("E38").Value = '  +4.46%  '

$ws.Range("D39").Value = "'2.791"
$ws.Range("E39").Value = '  -0.39%  '

$ws.Range("D40").Value = "'0.01806"
$ws.Range("E40").Value = '  +1.44%  '

$ws.Range("D41").Value = "'0.9044"
$ws.Range("E41").Value = '  -0.54%  '

$ws.Range("D42").Value = "'6.072"

$ws.Range("D43").Value = "'0.9991"
$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44").Value = "'101.37"
$ws.Range("E44").Value = '  -0.54%  '

$ws.Range("D45").Value = "'67.13"
$ws.Range("E45").Value = '  +1.55%  '

$ws.Range("D46").Value = "'7.315"
$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("D47").Value = "'9.160"
$ws.Range("E47").Value = '  +1.20%  '

$ws.Range("D48").Value = "'0.4009"
$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("D49").Value = "'1.695"
$ws.Range("E49").Value = '  +2.53%  '

$ws.Range("D50").Value = "'0.1122"
$ws.Range("E50").Value = '  -0.37%  '

$ws.Range("D51").Value = "'0.05743"
$ws.Range("E51").Value = '  -0.04%  '

# Reset number formatting on the Price column back to the default "Normal"
# style so the quote-prefix trick above does not leave any explicit cell
# style/format behind.
$ws.Range("D2:D51").Style = "Normal"
